# Applies the betexplorer scrape refresh described by the commit
# "Atualizado por script em 09-11-2023 02:45":
#   1) rows 192-195 (F:V, the match odds/urls) are cyclically rotated
#      up by one slot (row 192 data moves down into row 195) because
#      the upstream scrape re-sorted same-kickoff-time fixtures;
#   2) five brand-new fixtures are appended as rows 197-201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rotate the F:V match data across rows 192-195 -------------
# row 192
$ws.Cells.Item(192,6).Value = "Chico"
$ws.Cells.Item(192,7).Value = 1
$ws.Cells.Item(192,8).Value = "Dep. Cali"
$ws.Cells.Item(192,9).Value = 1
$ws.Cells.Item(192,10).Value = 2.65
$ws.Cells.Item(192,11).Value = "01/11/2023 10:43"
$ws.Cells.Item(192,12).Value = 3.75
$ws.Cells.Item(192,13).Value = "08/11/2023 01:23"
$ws.Cells.Item(192,14).Value = 2.97
$ws.Cells.Item(192,15).Value = "01/11/2023 10:43"
$ws.Cells.Item(192,16).Value = 3.1
$ws.Cells.Item(192,17).Value = "08/11/2023 01:23"
$ws.Cells.Item(192,18).Value = 2.9
$ws.Cells.Item(192,19).Value = "01/11/2023 10:43"
$ws.Cells.Item(192,20).Value = 2.25
$ws.Cells.Item(192,21).Value = "08/11/2023 01:23"
$ws.Cells.Item(192,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/chico-dep-cali/pGp3UN4o/"

# row 193
$ws.Cells.Item(193,6).Value = "Envigado"
$ws.Cells.Item(193,7).Value = 1
$ws.Cells.Item(193,8).Value = "Dep. Pasto"
$ws.Cells.Item(193,9).Value = 1
$ws.Cells.Item(193,10).Value = 2.7
$ws.Cells.Item(193,11).Value = "01/11/2023 10:42"
$ws.Cells.Item(193,12).Value = 2.83
$ws.Cells.Item(193,13).Value = "08/11/2023 01:29"
$ws.Cells.Item(193,14).Value = 2.92
$ws.Cells.Item(193,15).Value = "01/11/2023 10:42"
$ws.Cells.Item(193,16).Value = 3.26
$ws.Cells.Item(193,17).Value = "08/11/2023 01:29"
$ws.Cells.Item(193,18).Value = 3.03
$ws.Cells.Item(193,19).Value = "01/11/2023 10:42"
$ws.Cells.Item(193,20).Value = 2.69
$ws.Cells.Item(193,21).Value = "08/11/2023 01:29"
$ws.Cells.Item(193,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/envigado-dep-pasto/ptUg3o6n/"

# row 194
$ws.Cells.Item(194,6).Value = "Junior"
$ws.Cells.Item(194,7).Value = 2
$ws.Cells.Item(194,8).Value = "Huila"
$ws.Cells.Item(194,9).Value = 0
$ws.Cells.Item(194,10).Value = 1.51
$ws.Cells.Item(194,11).Value = "01/11/2023 10:42"
$ws.Cells.Item(194,12).Value = 1.3
$ws.Cells.Item(194,13).Value = "08/11/2023 01:28"
$ws.Cells.Item(194,14).Value = 4.13
$ws.Cells.Item(194,15).Value = "01/11/2023 10:42"
$ws.Cells.Item(194,16).Value = 5.33
$ws.Cells.Item(194,17).Value = "08/11/2023 01:28"
$ws.Cells.Item(194,18).Value = 7.03
$ws.Cells.Item(194,19).Value = "01/11/2023 10:42"
$ws.Cells.Item(194,20).Value = 12.68
$ws.Cells.Item(194,21).Value = "08/11/2023 01:28"
$ws.Cells.Item(194,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/junior-huila/8KEJcszO/"

# row 195
$ws.Cells.Item(195,6).Value = "Petrolera"
$ws.Cells.Item(195,7).Value = 2
$ws.Cells.Item(195,8).Value = "Pereira"
$ws.Cells.Item(195,9).Value = 1
$ws.Cells.Item(195,10).Value = 2.18
$ws.Cells.Item(195,11).Value = "01/11/2023 10:42"
$ws.Cells.Item(195,12).Value = 1.98
$ws.Cells.Item(195,13).Value = "08/11/2023 01:22"
$ws.Cells.Item(195,14).Value = 3.1
$ws.Cells.Item(195,15).Value = "01/11/2023 10:42"
$ws.Cells.Item(195,16).Value = 3.24
$ws.Cells.Item(195,17).Value = "08/11/2023 01:20"
$ws.Cells.Item(195,18).Value = 3.57
$ws.Cells.Item(195,19).Value = "01/11/2023 10:42"
$ws.Cells.Item(195,20).Value = 4.51
$ws.Cells.Item(195,21).Value = "08/11/2023 01:22"
$ws.Cells.Item(195,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/petrolera-dep-pereira/QePEbNKH/"

# --- Step 2: append 5 new fixtures as rows 197-201 ----------------------
# Clone row 196s formatting (bold/boxed index cell + datetime cell)
# onto each new row before filling in the values.
$ws.Range("A196:V196").Copy()
$ws.Range("A197:V197").PasteSpecial(-4122)
$ws.Range("A196:V196").Copy()
$ws.Range("A198:V198").PasteSpecial(-4122)
$ws.Range("A196:V196").Copy()
$ws.Range("A199:V199").PasteSpecial(-4122)
$ws.Range("A196:V196").Copy()
$ws.Range("A200:V200").PasteSpecial(-4122)
$ws.Range("A196:V196").Copy()
$ws.Range("A201:V201").PasteSpecial(-4122)

# row 197
$ws.Cells.Item(197,1).Value = 196
$ws.Cells.Item(197,2).Value = "colombia"
$ws.Cells.Item(197,3).Value = "primera-a"
$ws.Cells.Item(197,4).Value = "'2023"
$ws.Cells.Item(197,5).Value = 45239.0625
$ws.Cells.Item(197,6).Value = "America De Cali"
$ws.Cells.Item(197,7).Value = 1
$ws.Cells.Item(197,8).Value = "Bucaramanga"
$ws.Cells.Item(197,9).Value = 2
$ws.Cells.Item(197,10).Value = 1.46
$ws.Cells.Item(197,11).Value = "02/11/2023 01:42"
$ws.Cells.Item(197,12).Value = 1.4
$ws.Cells.Item(197,13).Value = "09/11/2023 01:26"
$ws.Cells.Item(197,14).Value = 4.38
$ws.Cells.Item(197,15).Value = "02/11/2023 01:42"
$ws.Cells.Item(197,16).Value = 4.66
$ws.Cells.Item(197,17).Value = "09/11/2023 01:29"
$ws.Cells.Item(197,18).Value = 7.32
$ws.Cells.Item(197,19).Value = "02/11/2023 01:42"
$ws.Cells.Item(197,20).Value = 8.09
$ws.Cells.Item(197,21).Value = "09/11/2023 01:29"
$ws.Cells.Item(197,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/america-de-cali-bucaramanga/IVL60qj5/"

# row 198
$ws.Cells.Item(198,1).Value = 197
$ws.Cells.Item(198,2).Value = "colombia"
$ws.Cells.Item(198,3).Value = "primera-a"
$ws.Cells.Item(198,4).Value = "'2023"
$ws.Cells.Item(198,5).Value = 45239.0625
$ws.Cells.Item(198,6).Value = "Atl. Nacional"
$ws.Cells.Item(198,7).Value = 2
$ws.Cells.Item(198,8).Value = "Deportes Tolima"
$ws.Cells.Item(198,9).Value = 3
$ws.Cells.Item(198,10).Value = 1.99
$ws.Cells.Item(198,11).Value = "02/11/2023 01:42"
$ws.Cells.Item(198,12).Value = 1.87
$ws.Cells.Item(198,13).Value = "09/11/2023 01:26"
$ws.Cells.Item(198,14).Value = 3.32
$ws.Cells.Item(198,15).Value = "02/11/2023 01:42"
$ws.Cells.Item(198,16).Value = 3.44
$ws.Cells.Item(198,17).Value = "09/11/2023 01:23"
$ws.Cells.Item(198,18).Value = 4.17
$ws.Cells.Item(198,19).Value = "02/11/2023 01:42"
$ws.Cells.Item(198,20).Value = 4.53
$ws.Cells.Item(198,21).Value = "09/11/2023 01:26"
$ws.Cells.Item(198,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/atl-nacional-deportes-tolima/MNINd1kU/"

# row 199
$ws.Cells.Item(199,1).Value = 198
$ws.Cells.Item(199,2).Value = "colombia"
$ws.Cells.Item(199,3).Value = "primera-a"
$ws.Cells.Item(199,4).Value = "'2023"
$ws.Cells.Item(199,5).Value = 45239.0625
$ws.Cells.Item(199,6).Value = "Jaguares de Cordoba"
$ws.Cells.Item(199,7).Value = 0
$ws.Cells.Item(199,8).Value = "Aguilas"
$ws.Cells.Item(199,9).Value = 1
$ws.Cells.Item(199,10).Value = 3.28
$ws.Cells.Item(199,11).Value = "02/11/2023 01:42"
$ws.Cells.Item(199,12).Value = 3.85
$ws.Cells.Item(199,13).Value = "09/11/2023 01:29"
$ws.Cells.Item(199,14).Value = 2.99
$ws.Cells.Item(199,15).Value = "02/11/2023 01:42"
$ws.Cells.Item(199,16).Value = 2.99
$ws.Cells.Item(199,17).Value = "09/11/2023 01:28"
$ws.Cells.Item(199,18).Value = 2.48
$ws.Cells.Item(199,19).Value = "02/11/2023 01:42"
$ws.Cells.Item(199,20).Value = 2.22
$ws.Cells.Item(199,21).Value = "09/11/2023 01:29"
$ws.Cells.Item(199,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/jaguares-de-cordoba-aguilas-doradas/dGXc25Lh/"

# row 200
$ws.Cells.Item(200,1).Value = 199
$ws.Cells.Item(200,2).Value = "colombia"
$ws.Cells.Item(200,3).Value = "primera-a"
$ws.Cells.Item(200,4).Value = "'2023"
$ws.Cells.Item(200,5).Value = 45239.0625
$ws.Cells.Item(200,6).Value = "La Equidad"
$ws.Cells.Item(200,7).Value = 2
$ws.Cells.Item(200,8).Value = "Millonarios"
$ws.Cells.Item(200,9).Value = 1
$ws.Cells.Item(200,10).Value = 2.39
$ws.Cells.Item(200,11).Value = "02/11/2023 01:42"
$ws.Cells.Item(200,12).Value = 2.11
$ws.Cells.Item(200,13).Value = "09/11/2023 01:26"
$ws.Cells.Item(200,14).Value = 2.94
$ws.Cells.Item(200,15).Value = "02/11/2023 01:42"
$ws.Cells.Item(200,16).Value = 3.02
$ws.Cells.Item(200,17).Value = "09/11/2023 01:26"
$ws.Cells.Item(200,18).Value = 3.52
$ws.Cells.Item(200,19).Value = "02/11/2023 01:42"
$ws.Cells.Item(200,20).Value = 4.16
$ws.Cells.Item(200,21).Value = "09/11/2023 01:26"
$ws.Cells.Item(200,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/la-equidad-millonarios/vuN21Pza/"

# row 201
$ws.Cells.Item(201,1).Value = 200
$ws.Cells.Item(201,2).Value = "colombia"
$ws.Cells.Item(201,3).Value = "primera-a"
$ws.Cells.Item(201,4).Value = "'2023"
$ws.Cells.Item(201,5).Value = 45239.0625
$ws.Cells.Item(201,6).Value = "U. Magdalena"
$ws.Cells.Item(201,7).Value = 0
$ws.Cells.Item(201,8).Value = "Ind. Medellin"
$ws.Cells.Item(201,9).Value = 4
$ws.Cells.Item(201,10).Value = 3.17
$ws.Cells.Item(201,11).Value = "02/11/2023 01:42"
$ws.Cells.Item(201,12).Value = 3.88
$ws.Cells.Item(201,13).Value = "09/11/2023 01:29"
$ws.Cells.Item(201,14).Value = 3.16
$ws.Cells.Item(201,15).Value = "02/11/2023 01:42"
$ws.Cells.Item(201,16).Value = 3.53
$ws.Cells.Item(201,17).Value = "09/11/2023 01:29"
$ws.Cells.Item(201,18).Value = 2.44
$ws.Cells.Item(201,19).Value = "02/11/2023 01:42"
$ws.Cells.Item(201,20).Value = 1.98
$ws.Cells.Item(201,21).Value = "09/11/2023 01:29"
$ws.Cells.Item(201,22).Value = "https://www.betexplorer.com/football/colombia/primera-a/union-magdalena-ind-medellin/bqfTHpLu/"

